$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataBF = New-Object "object[,]" 24,5
$dataBF[0,0] = 1.02
$dataBF[0,1] = 1.038730219718744
$dataBF[0,2] = 1.048421045966259
$dataBF[0,3] = 1.042383156219913
$dataBF[0,4] = 1.057002457729889
$dataBF[1,0] = 1.02
$dataBF[1,1] = 1.039594736634425
$dataBF[1,2] = 1.049122234403329
$dataBF[1,3] = 1.04319603097846
$dataBF[1,4] = 1.05783466905203
$dataBF[2,0] = 1.02
$dataBF[2,1] = 1.040154695019364
$dataBF[2,2] = 1.049576399328184
$dataBF[2,3] = 1.043722902999394
$dataBF[2,4] = 1.058373954089918
$dataBF[3,0] = 1.02
$dataBF[3,1] = 1.040390233688805
$dataBF[3,2] = 1.049767436118342
$dataBF[3,3] = 1.043944610810666
$dataBF[3,4] = 1.058600856545493
$dataBF[4,0] = 1.02
$dataBF[4,1] = 1.040429789375829
$dataBF[4,2] = 1.049799518195043
$dataBF[4,3] = 1.043981848840457
$dataBF[4,4] = 1.058638965381669
$dataBF[5,0] = 1.02
$dataBF[5,1] = 1.040157841781091
$dataBF[5,2] = 1.049578951557274
$dataBF[5,3] = 1.043725864643802
$dataBF[5,4] = 1.058376985238988
$dataBF[6,0] = 1.02
$dataBF[6,1] = 1.039022270818278
$dataBF[6,2] = 1.04865792188817
$dataBF[6,3] = 1.04265768622318
$dataBF[6,4] = 1.057283543073003
$dataBF[7,0] = 1.02
$dataBF[7,1] = 1.037025588640066
$dataBF[7,2] = 1.047038463184446
$dataBF[7,3] = 1.040782293341072
$dataBF[7,4] = 1.055362877278547
$dataBF[8,0] = 1.02
$dataBF[8,1] = 1.035697470677548
$dataBF[8,2] = 1.045961285176571
$dataBF[8,3] = 1.039536754411765
$dataBF[8,4] = 1.054086653530341
$dataBF[9,0] = 1.02
$dataBF[9,1] = 1.03512311134586
$dataBF[9,2] = 1.04549545980137
$dataBF[9,3] = 1.038998563857334
$dataBF[9,4] = 1.053535057825792
$dataBF[10,0] = 1.02
$dataBF[10,1] = 1.034909878997121
$dataBF[10,2] = 1.045322523308764
$dataBF[10,3] = 1.038798828314865
$dataBF[10,4] = 1.053330325300452
$dataBF[11,0] = 1.02
$dataBF[11,1] = 1.034955613046325
$dataBF[11,2] = 1.045359614592738
$dataBF[11,3] = 1.038841664438041
$dataBF[11,4] = 1.053374234093724
$dataBF[12,0] = 1.02
$dataBF[12,1] = 1.035105483235419
$dataBF[12,2] = 1.045481162934702
$dataBF[12,3] = 1.038982050128399
$dataBF[12,4] = 1.053518131405587
$dataBF[13,0] = 1.02
$dataBF[13,1] = 1.035197837861356
$dataBF[13,2] = 1.045556065093372
$dataBF[13,3] = 1.039068569280594
$dataBF[13,4] = 1.05360681183885
$dataBF[14,0] = 1.02
$dataBF[14,1] = 1.035735604401436
$dataBF[14,2] = 1.045992213253129
$dataBF[14,3] = 1.039572496427686
$dataBF[14,4] = 1.054123282713604
$dataBF[15,0] = 1.02
$dataBF[15,1] = 1.036073126125329
$dataBF[15,2] = 1.046265959241982
$dataBF[15,3] = 1.039888901971454
$dataBF[15,4] = 1.054447524956838
$dataBF[16,0] = 1.02
$dataBF[16,1] = 1.036270066645036
$dataBF[16,2] = 1.046425688406422
$dataBF[16,3] = 1.040073565466345
$dataBF[16,4] = 1.054636748034785
$dataBF[17,0] = 1.02
$dataBF[17,1] = 1.036337230070125
$dataBF[17,2] = 1.046480161669854
$dataBF[17,3] = 1.040136549465846
$dataBF[17,4] = 1.054701284792977
$dataBF[18,0] = 1.02
$dataBF[18,1] = 1.036036906007663
$dataBF[18,2] = 1.046236582901525
$dataBF[18,3] = 1.039854943313121
$dataBF[18,4] = 1.054412726676632
$dataBF[19,0] = 1.02
$dataBF[19,1] = 1.035061347144356
$dataBF[19,2] = 1.045445367430111
$dataBF[19,3] = 1.038940705268938
$dataBF[19,4] = 1.053475752949462
$dataBF[20,0] = 1.02
$dataBF[20,1] = 1.0344486133268
$dataBF[20,2] = 1.044948430773555
$dataBF[20,3] = 1.038366885784161
$dataBF[20,4] = 1.052887536193707
$dataBF[21,0] = 1.02
$dataBF[21,1] = 1.034773374090908
$dataBF[21,2] = 1.045211815318302
$dataBF[21,3] = 1.038670983152909
$dataBF[21,4] = 1.053199275512675
$dataBF[22,0] = 1.02
$dataBF[22,1] = 1.036053272104218
$dataBF[22,2] = 1.046249856627023
$dataBF[22,3] = 1.039870287430978
$dataBF[22,4] = 1.054428450218801
$dataBF[23,0] = 1.02
$dataBF[23,1] = 1.03754125578448
$dataBF[23,2] = 1.047456705532727
$dataBF[23,3] = 1.041266302312811
$dataBF[23,4] = 1.055858679364316
$ws.Range("B2:F25").Value2 = $dataBF

$dataIN = New-Object "object[,]" 24,6
$dataIN[0,0] = 1.044812931410447
$dataIN[0,1] = 1.04382607270359
$dataIN[0,2] = 1.051180938695051
$dataIN[0,3] = 1.045159998979641
$dataIN[0,4] = 1.059738621068969
$dataIN[0,5] = 1.018523846863181
$dataIN[1,0] = 1.045066039929225
$dataIN[1,1] = 1.044335929087231
$dataIN[1,2] = 1.051694456001866
$dataIN[1,3] = 1.045783688715999
$dataIN[1,4] = 1.060384535801566
$dataIN[1,5] = 1.018694773103118
$dataIN[2,0] = 1.045228803824759
$dataIN[2,1] = 1.044665739716124
$dataIN[2,2] = 1.05202650868087
$dataIN[2,3] = 1.046187492884141
$dataIN[2,4] = 1.060802623703341
$dataIN[2,5] = 1.018805293783896
$dataIN[3,0] = 1.045296986181913
$dataIN[3,1] = 1.044804366864484
$dataIN[3,2] = 1.052166047942837
$dataIN[3,3] = 1.046357307084524
$dataIN[3,4] = 1.060978419385493
$dataIN[3,5] = 1.018851737074895
$dataIN[4,0] = 1.045308419993343
$dataIN[4,1] = 1.044827641458612
$dataIN[4,2] = 1.052189473897104
$dataIN[4,3] = 1.046385822816884
$dataIN[4,4] = 1.061007938059871
$dataIN[4,5] = 1.018859533945338
$dataIN[5,0] = 1.045229715838836
$dataIN[5,1] = 1.044667592159506
$dataIN[5,2] = 1.052028373431691
$dataIN[5,3] = 1.046189761735698
$dataIN[5,4] = 1.060804972571822
$dataIN[5,5] = 1.018805914439008
$dataIN[6,0] = 1.04489868015936
$dataIN[6,1] = 1.043998401170663
$dataIN[6,2] = 1.051354530594342
$dataIN[6,3] = 1.045370728315729
$dataIN[6,4] = 1.059956881538276
$dataIN[6,5] = 1.01858162848875
$dataIN[7,0] = 1.04430761934233
$dataIN[7,1] = 1.042818482798831
$dataIN[7,2] = 1.050165447604491
$dataIN[7,3] = 1.043929349179006
$dataIN[7,4] = 1.058463558199866
$dataIN[7,5] = 1.018185814738878
$dataIN[8,0] = 1.043908425104648
$dataIN[8,1] = 1.042031457071572
$dataIN[8,2] = 1.049371666753223
$dataIN[8,3] = 1.042969763990054
$dataIN[8,4] = 1.057468851315885
$dataIN[8,5] = 1.017921566482151
$dataIN[9,0] = 1.043734356834505
$dataIN[9,1] = 1.041690582756976
$dataIN[9,2] = 1.049027715399291
$dataIN[9,3] = 1.042554586716932
$dataIN[9,4] = 1.057038350468638
$dataIN[9,5] = 1.017807061293947
$dataIN[10,0] = 1.043669518412261
$dataIN[10,1] = 1.041563954953489
$dataIN[10,2] = 1.04889992189254
$dataIN[10,3] = 1.042400422268333
$dataIN[10,4] = 1.05687847688761
$dataIN[10,5] = 1.017764516800208
$dataIN[11,0] = 1.043683434684863
$dataIN[11,1] = 1.041591117567269
$dataIN[11,2] = 1.048927335581903
$dataIN[11,3] = 1.042433488741441
$dataIN[11,4] = 1.056912768769114
$dataIN[11,5] = 1.017773643282059
$dataIN[12,0] = 1.043729000972889
$dataIN[12,1] = 1.041680115899384
$dataIN[12,2] = 1.049017152647977
$dataIN[12,3] = 1.042541842389668
$dataIN[12,4] = 1.057025134574527
$dataIN[12,5] = 1.017803544798489
$dataIN[13,0] = 1.043757051810214
$dataIN[13,1] = 1.041734949166289
$dataIN[13,2] = 1.049072487351369
$dataIN[13,3] = 1.042608609431395
$dataIN[13,4] = 1.057094371352642
$dataIN[13,5] = 1.01782196651047
$dataIN[14,0] = 1.043919951917929
$dataIN[14,1] = 1.042054078043792
$dataIN[14,2] = 1.04939448874277
$dataIN[14,3] = 1.042997324997224
$dataIN[14,4] = 1.057497426870172
$dataIN[14,5] = 1.01792916408356
$dataIN[15,0] = 1.044021810218364
$dataIN[15,1] = 1.042254236690324
$dataIN[15,2] = 1.049596408600714
$dataIN[15,3] = 1.043241245178633
$dataIN[15,4] = 1.05775031123278
$dataIN[15,5] = 1.017996384126866
$dataIN[16,0] = 1.044081105206447
$dataIN[16,1] = 1.042370977372023
$dataIN[16,2] = 1.049714161893143
$dataIN[16,3] = 1.043383551387776
$dataIN[16,4] = 1.057897834932048
$dataIN[16,5] = 1.018035584311204
$dataIN[17,0] = 1.044101303357046
$dataIN[17,1] = 1.042410781455555
$dataIN[17,2] = 1.049754308759156
$dataIN[17,3] = 1.043432079461443
$dataIN[17,4] = 1.057948140138171
$dataIN[17,5] = 1.018048949168602
$dataIN[18,0] = 1.044010893907459
$dataIN[18,1] = 1.042232762429574
$dataIN[18,2] = 1.049574746895834
$dataIN[18,3] = 1.043215071563361
$dataIN[18,4] = 1.05772317699984
$dataIN[18,5] = 1.017989172890066
$dataIN[19,0] = 1.043715587846975
$dataIN[19,1] = 1.041653908428115
$dataIN[19,2] = 1.048990704707389
$dataIN[19,3] = 1.042509933525378
$dataIN[19,4] = 1.056992044706009
$dataIN[19,5] = 1.017794739879522
$dataIN[20,0] = 1.043528866225937
$dataIN[20,1] = 1.041289891636521
$dataIN[20,2] = 1.048623294292529
$dataIN[20,3] = 1.042066880293207
$dataIN[20,4] = 1.056532547956825
$dataIN[20,5] = 1.01767242194526
$dataIN[21,0] = 1.043627950258826
$dataIN[21,1] = 1.04148286991147
$dataIN[21,2] = 1.048818084037166
$dataIN[21,3] = 1.042301722766882
$dataIN[21,4] = 1.056776116885613
$dataIN[21,5] = 1.017737271500928
$dataIN[22,0] = 1.04401582688079
$dataIN[22,1] = 1.042242465749964
$dataIN[22,2] = 1.04958453495951
$dataIN[22,3] = 1.043226898195705
$dataIN[22,4] = 1.05773543772847
$dataIN[22,5] = 1.017992431362744
$dataIN[23,0] = 1.044461334478744
$dataIN[23,1] = 1.043123597950155
$dataIN[23,2] = 1.050165447604491
$dataIN[23,3] = 1.043929349179006
$dataIN[23,4] = 1.058849476114547
$dataIN[23,5] = 1.018288209658643
$ws.Range("I2:N25").Value2 = $dataIN
